$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 4 (pushes old rows 4..29 down to 6..31).
$ws.Rows.Item(4).Resize(2).Insert()

# The insert leaves the new A4:A5 cells unstyled; pull the bordered/bold
# style used by the rest of column A (row 6, formerly row 4) without
# disturbing the (currently empty) values.
$ws.Range("A6").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New HKL entry "Holden" (index 2) with freshly simulated data.
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 0.8833391104341288
$ws.Range("D4").Value = 1.679608811178225
$ws.Range("E4").Value = 0.8440207403301154
$ws.Range("F4").Value = 1.679608811178225
$ws.Range("G4").Value = 0.8440207403301154
$ws.Range("H4").Value = 0.8091354226605793
$ws.Range("I4").Value = 1.235638612406367
$ws.Range("J4").Value = 0.8717882282687709
$ws.Range("K4").Value = 0.8440207403301154
$ws.Range("L4").Value = 0.8833391104341288
$ws.Range("M4").Value = 1.281473960806177
$ws.Range("N4").Value = 1.281473960806177
$ws.Range("O4").Value = 1.266195511339574
$ws.Range("P4").Value = 1.13565622064749
$ws.Range("Q4").Value = 1.13565622064749
$ws.Range("R4").Value = 1.062747350568146
$ws.Range("S4").Value = 1.062747350568146
$ws.Range("T4").Value = 1.053921820879698

# New HKL entry "Rizzie Spiral" (index 3) with freshly simulated data.
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.188537314823157
$ws.Range("D5").Value = 0.9418479720296867
$ws.Range("E5").Value = 0.8152802329331432
$ws.Range("F5").Value = 0.9418479720296867
$ws.Range("G5").Value = 0.8152802329331432
$ws.Range("H5").Value = 1.918793037201898
$ws.Range("I5").Value = 0.724705281272511
$ws.Range("J5").Value = 1.048841773471372
$ws.Range("K5").Value = 0.8152802329331432
$ws.Range("L5").Value = 1.188537314823157
$ws.Range("M5").Value = 1.065192643426422
$ws.Range("N5").Value = 1.065192643426422
$ws.Range("O5").Value = 0.9516968560417851
$ws.Range("P5").Value = 0.9818885065953292
$ws.Range("Q5").Value = 0.9818885065953292
$ws.Range("R5").Value = 0.9402364381797828
$ws.Range("S5").Value = 0.9402364381797828
$ws.Range("T5").Value = 1.106334268621961

# Rename "Thomas Hex" -> "Matthies Hex" (the row that shifted from 9 to 11).
$ws.Range("B11").Value = "Matthies Hex"
